# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.831.02"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "1.617.48"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("D4").Value = "'0.993"
$ws.Range("E4").Value = "  -0.79%  "
$ws.Range("D5").Value = "'213.24"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").Value = "  -0.86%  "
$ws.Range("D8").Value = "'29.12"
$ws.Range("E8").Value = "  +8.89%  "
$ws.Range("E9").Value = "  +3.25%  "
$ws.Range("D10").Value = "'0.0605"
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "1.850.06"
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").Value = "1.622.46"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("E14").Value = "  +5.68%  "
$ws.Range("E15").Value = "  +4.97%  "
$ws.Range("D16").Value = "29.861.44"
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("E17").Value = "  +15.68%  "
$ws.Range("D18").Value = "'64.34"
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("D19").Value = "'240.86"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("E20").Value = "  +2.45%  "
$ws.Range("D21").Value = "'0.994"
$ws.Range("E21").Value = "  -0.66%  "
$ws.Range("E22").Value = "  +2.48%  "
$ws.Range("D23").Value = "'9.57"
$ws.Range("E23").Value = "  +4.06%  "
$ws.Range("D24").Value = "'2.12"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("D25").Value = "'155.15"
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("D26").Value = "'15.58"
$ws.Range("E26").Value = "  +2.10%  "
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("E28").Value = "  +2.96%  "
$ws.Range("D29").Value = "'0.994"
$ws.Range("E29").Value = "  -0.72%  "
$ws.Range("E30").Value = "  +3.07%  "
$ws.Range("E31").Value = "  +5.48%  "
$ws.Range("D33").Value = "'3.20"
$ws.Range("E33").Value = "  +3.53%  "
$ws.Range("D34").Value = "1.415.83"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("E35").Value = "  +6.27%  "
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("D37").Value = "'2.87"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("D38").Value = "'2.29"
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("E39").Value = "  +2.32%  "
$ws.Range("D40").Value = "'0.556"
$ws.Range("E40").Value = "  +3.34%  "
$ws.Range("D41").Value = "'0.0503"
$ws.Range("E41").Value = "  +3.03%  "
$ws.Range("E42").Value = "  +3.53%  "
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'69.12"
$ws.Range("E44").Value = "  +5.12%  "
$ws.Range("B45").Value = "BitcoinSV"
$ws.Range("C45").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D45").Value = "'53.37"
$ws.Range("E45").Value = "  +1.21%  "
$ws.Range("E46").Value = "  +18.98%  "
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("E48").Value = "  +2.92%  "
$ws.Range("D49").Value = "1.759.05"
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("D50").Value = "'88.11"
$ws.Range("E50").Value = "  +1.66%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0107"
$ws.Range("E51").Value = "  +4.09%  "
